$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "r775"
$ws.Range("B11").Value = "bruce"
$ws.Range("C11").Value = "we might be close to ready "
$ws.Range("D11").Value = "2025-10-01 14:44:36"
